$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The cell V2 held '<mods:note displayLabel="Description">' but with a
# non-breaking space (U+00A0) between "mods:note" and "displayLabel"
# instead of a normal space. Rewrite it with a normal ASCII space.
$ws.Range("V2").Value = '<mods:note displayLabel="Description">'
